$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated AC cable incidence matrix (B2:I9)
$data = @(
    @(0,1,1,0,1,1,0,1),
    @(1,0,1,0,0,1,1,0),
    @(1,1,0,0,0,0,0,1),
    @(0,0,0,0,0,1,1,0),
    @(1,0,0,0,0,1,0,0),
    @(1,1,0,1,1,0,0,0),
    @(0,1,0,1,0,0,0,0),
    @(1,0,1,0,0,0,0,0)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = 2 + $i
    $rowData = $data[$i]
    for ($j = 0; $j -lt $rowData.Length; $j++) {
        $col = 2 + $j
        $ws.Cells.Item($row, $col).Value = $rowData[$j]
    }
}

$wb.Save()
